$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 1.41
$ws.Range("H2").Value = 10
$ws.Range("N2").Value = 3.5
$ws.Range("O2").Value = 1.37
$ws.Range("Q2").Value = 2.1
$ws.Range("S2").Value = 3.95
$ws.Range("X2").Value = 1000
$ws.Range("Y2").Value = 27
$ws.Range("AC2").Value = 12
$ws.Range("AI2").Value = 270
$ws.Range("AJ2").Value = 980
$ws.Range("AM2").Value = 1000
$ws.Range("O3").Value = 1.32
$ws.Range("Q3").Value = 1.99
$ws.Range("T3").Value = 1.76
$ws.Range("X3").Value = 1000
$ws.Range("Y3").Value = 12
$ws.Range("Z3").Value = 1000
$ws.Range("AA3").Value = 1000
$ws.Range("AB3").Value = 13.5
$ws.Range("AD3").Value = 14
$ws.Range("AE3").Value = 1000
$ws.Range("AF3").Value = 1000
$ws.Range("AG3").Value = 1000
$ws.Range("AH3").Value = 1000
$ws.Range("AI3").Value = 1000
$ws.Range("AJ3").Value = 1000
$ws.Range("AK3").Value = 1000
$ws.Range("AL3").Value = 1000
$ws.Range("AM3").Value = 1000
$ws.Range("AN3").Value = 1000
$ws.Range("AO3").Value = 1000
$ws.Range("I4").Value = 8.6
$ws.Range("N4").Value = 7.2
$ws.Range("O4").Value = 1.14
$ws.Range("Q4").Value = 1.44
$ws.Range("S4").Value = 2.12
$ws.Range("T4").Value = 1.68
$ws.Range("U4").Value = 2.38
$ws.Range("X4").Value = 1000
$ws.Range("Y4").Value = 1000
$ws.Range("Z4").Value = 1000
$ws.Range("AA4").Value = 1000
$ws.Range("AB4").Value = 980
$ws.Range("AD4").Value = 1000
$ws.Range("AE4").Value = 100
$ws.Range("AI4").Value = 1000
$ws.Range("AM4").Value = 90
$ws.Range("AO4").Value = 1000
$ws.Range("G5").Value = 2.2
$ws.Range("N5").Value = 4.6
$ws.Range("R5").Value = 1.48
$ws.Range("X5").Value = 18.5
$ws.Range("Y5").Value = 16
$ws.Range("Z5").Value = 26
$ws.Range("AA5").Value = 1000
$ws.Range("AD5").Value = 15
$ws.Range("AI5").Value = 55
$ws.Range("AN5").Value = 14
$ws.Range("AO5").Value = 30
$ws.Range("F6").Value = 1.85
$ws.Range("H6").Value = 5.2
$ws.Range("I6").Value = 5.4
$ws.Range("K6").Value = 3.7
$ws.Range("M6").Value = 1.09
$ws.Range("P6").Value = 1.81
$ws.Range("U6").Value = 1.91
$ws.Range("Z6").Value = 40
$ws.Range("AE6").Value = 80
$ws.Range("T7").Value = 1.77
$ws.Range("AE7").Value = 27
$ws.Range("AF7").Value = 23
$ws.Range("AN7").Value = 36
$ws.Range("K8").Value = 3.8
$ws.Range("U8").Value = 2
$ws.Range("AH8").Value = 21
$ws.Range("N9").Value = 4.3
$ws.Range("Q9").Value = 1.85
$ws.Range("R9").Value = 1.44
$ws.Range("T9").Value = 1.73
$ws.Range("Z9").Value = 15
$ws.Range("AC9").Value = 8.4
$ws.Range("AH9").Value = 16.5
$ws.Range("AN9").Value = 40
$ws.Range("K10").Value = 3.3
$ws.Range("F11").Value = 2.36
$ws.Range("AH11").Value = 21
$ws.Range("AN11").Value = 29
$ws.Range("I12").Value = 1.45
$ws.Range("O12").Value = 1.22
$ws.Range("P12").Value = 2.4
$ws.Range("S12").Value = 2.7
$ws.Range("G13").Value = 1.75
$ws.Range("P13").Value = 2.08
$ws.Range("Z13").Value = 50
$ws.Range("AF13").Value = 10.5
$ws.Range("Q14").Value = 1.81
